# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Delete the two rows that belong to the old "2503" period (bottom-most first so row
# numbers of rows above are not affected by the deletion).
$ws.Rows.Item(27).Delete()   # CC / 7918671 / ARIEL PUELLO MONTERROSA / 2503
$ws.Rows.Item(25).Delete()   # CC / 1235048469 / YELITZA MARIA CAÑATE POLO / 2503

# The remaining "2504" rows (now at 24 and 25 after the deletions above) become the new
# "2509" period.
$ws.Range("E24").Value = "2509"
$ws.Range("E25").Value = "2509"

# Update the summary fields: total overdue value and period count.
$ws.Range("E11").Value = 569300
$ws.Range("F13").Value = 5

# The last row of the data table carries a solid (non-auto) black bottom border to close
# the table visually. Row 25 is now the last data row, so give it that border.
$ws.Range("B25:J25").Borders.Item(9).Color = 0
